$d = $word.ActiveDocument

$xml20 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Agile software development agency. </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Full or</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> part-time building solutions across many industries.</w:t></w:r><w:r><w:rPr><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2016-18 </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Contract w/ </w:t></w:r><w:hyperlink r:id="rId42"><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Intricately</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">authoritative product adoption, usage, and spend Data on 7 million companies with m</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">arketing and sales intelligence.</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> Built systems to track changes to </w:t></w:r><w:r><w:rPr><w:i w:val="1"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">the internet</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> (writ large). Automated alerting, and benchmarking. Cloud intel r</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">eports for F100 enterprise customers. Wrote </w:t></w:r><w:hyperlink r:id="rId43"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">debug_logging</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">; maintaining </w:t></w:r><w:hyperlink r:id="rId44"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">dynamoid</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2018-19 Contract /w </w:t></w:r><w:hyperlink r:id="rId45"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Group1001</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">family of insurance companies offering accumulation and protection solutions;</w:t></w:r><w:r><w:rPr><w:color w:val="3c4858"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="f9fafc" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId46"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Gainbridge</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:hyperlink r:id="rId47"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Delaware Life</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">, Relay Rewards</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">. Sole dev on Relay Rewards &amp; Group1001 (Rails); Analytics SME on Gainbridge (React v16)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml49 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2013-14 Contract w/ </w:t></w:r><w:hyperlink r:id="rId75"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Church Pension Group</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">;</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Provided resources and services to sustain the operations of the Episcopal Church worldwide</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:fill="auto" w:val="clear"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">. Forensic modernization; Implemented  token auth, daemons, `make` system. 10x perf improvement of internal servi</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">ces, 20 req/s increased to &gt; 200 req/s</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r20 = $d.Paragraphs.Item(20).Range
$r20.MoveEnd(1, -1) | Out-Null
$r20 = $d.Range($r20.Start, $r20.End)
$r20.InsertXML($xml20)

$r49 = $d.Paragraphs.Item(49).Range
$r49.MoveEnd(1, -1) | Out-Null
$r49 = $d.Range($r49.Start, $r49.End)
$r49.InsertXML($xml49)

Write-Output "DONE"
